# Adds two new columns, I ("I0") and J ("IF"), to the sheet -
# mirrors the header style of the existing last header cell (H1)
# and fills in the per-row values for rows 2..64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the rest of the header row (H1 etc.):
# bold font, centered/top aligned, thin box border.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("I1:J1").VerticalAlignment = -4160    # xlTop
$ws.Range("I1:J1").Borders.LineStyle = 1        # xlContinuous

# Per-row (I, J) values for data rows 2..64.
$values = @(
    @(7,7),
    @(6,6),
    @(9,9),
    @(7,7),
    @(5,6),
    @(6,6),
    @(4,4),
    @(9,9),
    @(5,5),
    @(5,5),
    @(6,6),
    @(7,7),
    @(5,5),
    @(8,8),
    @(6,6),
    @(7,7),
    @(8,8),
    @(5,5),
    @(5,5),
    @(6,6),
    @(9,9),
    @(2,3),
    @(6,6),
    @(8,8),
    @(7,7),
    @(9,9),
    @(6,6),
    @(8,8),
    @(8,8),
    @(9,9),
    @(2,3),
    @(6,6),
    @(9,9),
    @(5,6),
    @(5,5),
    @(9,9),
    @(9,9),
    @(9,9),
    @(7,8),
    @(8,8),
    @(5,5),
    @(9,9),
    @(9,9),
    @(7,7),
    @(3,3),
    @(7,7),
    @(8,8),
    @(10,10),
    @(9,9),
    @(9,10),
    @(8,8),
    @(9,9),
    @(9,9),
    @(7,8),
    @(5,5),
    @(4,5),
    @(8,8),
    @(8,8),
    @(5,6),
    @(4,4),
    @(8,8),
    @(5,5),
    @(8,8)
)

for ($idx = 0; $idx -lt $values.Count; $idx++) {
    $row = $idx + 2
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value  = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]   # column J
}
